$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B10").Value = "hlthst_duration_cont_log_scale_clst"
$ws.Range("B3").Select()
